$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column C
$ws.Range("C2").Value = 11
$ws.Range("C3").Value = 9.5

# Update the active selection to C3
$ws.Range("C3").Select()
